$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16 - this shifts the existing rows 16-31 down to 17-32,
# matching the data already present in those rows (no other edits needed for them).
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly record.
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = "Vega Modelo de Temuco"
$ws.Range("C16").Value = "La Araucanía"
$ws.Range("D16").Value = 45079
$ws.Range("E16").Value = 9
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100108
$ws.Range("H16").Value = "Tropicales y subtropicales"
$ws.Range("I16").Value = 100108001
$ws.Range("J16").Value = "Guayaba"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 2600
$ws.Range("O16").Value = 2600
$ws.Range("P16").Value = 2600
$ws.Range("Q16").Value = "$/kilo"
$ws.Range("R16").Value = "Región de Arica y Parinacota"
$ws.Range("S16").Value = 2600
$ws.Range("T16").Value = 1

# Match the date-number format used by the other rows in column D.
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
